# Actualizacion automatica del tracker: agrega los nuevos resultados (filas 79-87)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-TrackerRow {
    param($Row, $EventId, $Fecha, $JugadorA, $JugadorB, $Pronostico, $Cuota)

    $ws.Cells.Item($Row, 1).Value = $EventId

    # Forzar texto en la columna de fecha para que Excel no la convierta
    # automaticamente en un numero de serie de fecha.
    $ws.Cells.Item($Row, 2).NumberFormat = "@"
    $ws.Cells.Item($Row, 2).Value = $Fecha
    $ws.Cells.Item($Row, 2).Style = "Normal"

    $ws.Cells.Item($Row, 3).Value = $JugadorA
    $ws.Cells.Item($Row, 4).Value = $JugadorB
    $ws.Cells.Item($Row, 5).Value = $Pronostico
    $ws.Cells.Item($Row, 6).Value = $Cuota
}

Add-TrackerRow 79 14494935 "2025-09-03" "Novak Djokovic"         "Taylor Fritz"          "Gana Taylor Fritz"              2.3
Add-TrackerRow 80 14592780 "2025-09-02" "Alex Martinez"          "Carlos Taberner"       "Gana Alex Martinez"             6
Add-TrackerRow 81 14592781 "2025-09-02" "Carlos Lopez Montagud"  "Carlos Sanchez Jover"  "Gana Carlos Lopez Montagud"     2.25
Add-TrackerRow 82 14592778 "2025-09-02" "Guy Den Ouden"          "Dominik Kellovsky"     "Gana Dominik Kellovsky"         6
Add-TrackerRow 83 14552613 "2025-09-03" "Sumit Nagal"            "Filip Cristian Jianu"  "Gana Filip Cristian Jianu"      3
Add-TrackerRow 84 14585167 "2025-09-02" "Justin Engel"           "Felix Corwin"          "Gana Felix Corwin"              3.25
Add-TrackerRow 85 14585168 "2025-09-02" "Nicolas Tepmahc"        "Maxime Janvier"        "Gana Nicolas Tepmahc"           2.62
Add-TrackerRow 86 14585371 "2025-09-02" "Hamish Stewart"         "Egor Gerasimov"        "Gana Egor Gerasimov"            3.4
Add-TrackerRow 87 14585908 "2025-09-02" "Jakub Paul"             "David Pichler"         "Gana David Pichler"             5
